# Updates the cryptos list (Price column D, Volume(1h) column E) with
# freshly scraped values, as produced by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text format on the touched cells so Excel keeps the exact
# textual representation (e.g. trailing zeros, thousands separators,
# and the padded "  +x.xx%  " volume strings) instead of coercing them
# into numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.428.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.649.58"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.35"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.02"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.83%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.648.97"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.04"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.132.99"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.314.97"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.634.64"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.56"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "363.77"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.41"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.17%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.64"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.88"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.780.76"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "565.31"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.08"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.40"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.66"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.89%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.58"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.64"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.370"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.87"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.34"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.65"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.29%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "158.30"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.83"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.99"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.73%  "
